$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.62
$ws.Range("I2").Value = 2.92
$ws.Range("L2").Value = 1.41
$ws.Range("Q2").Value = 1.96
$ws.Range("V2").Value = 1.52

# Row 3
$ws.Range("X3").Value = 32
$ws.Range("Y3").Value = 19
$ws.Range("Z3").Value = 21
$ws.Range("AB3").Value = 25
$ws.Range("AD3").Value = 13.5
$ws.Range("AF3").Value = 34
$ws.Range("AG3").Value = 17.5
$ws.Range("AH3").Value = 14.5
$ws.Range("AI3").Value = 28
$ws.Range("AK3").Value = 34
$ws.Range("AL3").Value = 36
$ws.Range("AN3").Value = 18.5

# Row 6
$ws.Range("S6").Value = 2.86
$ws.Range("T6").Value = 1.67

# Row 7
$ws.Range("AH7").Value = 24
$ws.Range("AJ7").Value = 21

# Row 8
$ws.Range("T8").Value = 2.36

# Row 9
$ws.Range("H9").Value = 2.16
$ws.Range("J9").Value = 3.3

# Row 10
$ws.Range("H10").Value = 3.5
$ws.Range("K10").Value = 4.1
